$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'schubert-winterreise_120'
$ws.Range("B2").Value = 'schubert-winterreise_14'
$ws.Range("C2").Value = 0.2196969696969697
$ws.Range("D2").Value = '[[''A#:min'', ''F:maj/A'', ''A#:min'']]'
$ws.Range("E2").Value = '[[''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("F2").Value = '[(15.32, 17.86)]'
$ws.Range("G2").Value = '[(10.56, 17.34)]'
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""

# Row 3
$ws.Range("A3").Value = 'isophonics_215'
$ws.Range("B3").Value = 'isophonics_19'
$ws.Range("C3").Value = 0.1416666666666667
$ws.Range("D3").Value = '[[''F:7'', ''Bb'', ''C'']]'
$ws.Range("E3").Value = '[[''F:7'', ''Bb'', ''C'']]'
$ws.Range("F3").Value = '[(2.572176, 6.09)]'
$ws.Range("G3").Value = '[(9.297913, 12.780907)]'
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

# Row 4
$ws.Range("A4").Value = 'isophonics_35'
$ws.Range("B4").Value = 'schubert-winterreise_214'
$ws.Range("C4").Value = 0.2266666666666667
$ws.Range("D4").Value = '[[''G:7'', ''C'', ''E'']]'
$ws.Range("E4").Value = '[[''D:7'', ''G:maj'', ''B:maj/D#'']]'
$ws.Range("F4").Value = '[(20.769931, 26.946439)]'
$ws.Range("G4").Value = '[(76.32, 84.1)]'
$ws.Range("H4").Value = 'spotify:track:1dxbAIfCASqv6jix2R1Taj'
$ws.Range("I4").Value = ""

# Row 5
$ws.Range("A5").Value = 'jaah_39'
$ws.Range("B5").Value = 'jaah_43'
$ws.Range("C5").Value = 0.09377038486627527
$ws.Range("D5").Value = '[[''Eb:7'', ''Eb:7'', ''Ab'', ''Ab'', ''Db'']]'
$ws.Range("E5").Value = '[[''Bb:7'', ''Bb:7'', ''Eb'', ''Eb'', ''Ab'']]'
$ws.Range("F5").Value = '[(121.03, 126.19)]'
$ws.Range("G5").Value = '[(93.27, 102.97)]'
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("A6").Value = 'schubert-winterreise_143'
$ws.Range("B6").Value = 'schubert-winterreise_164'
$ws.Range("C6").Value = 0.1846153846153846
$ws.Range("D6").Value = '[[''F#:7/B'', ''B:min'', ''E:min/B'', ''B:min'']]'
$ws.Range("E6").Value = '[[''E:7'', ''A:min'', ''D:min'', ''A:min'']]'
$ws.Range("F6").Value = '[(23.14, 32.76)]'
$ws.Range("G6").Value = '[(51.62, 63.9)]'
$ws.Range("H6").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Range("I6").Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'

# Row 7
$ws.Range("A7").Value = 'isophonics_6'
$ws.Range("B7").Value = 'jaah_32'
$ws.Range("C7").Value = 0.2209443099273608
$ws.Range("D7").Value = '[[''Bb'', ''Eb:7'', ''Bb'', ''Eb:7'', ''Bb'']]'
$ws.Range("E7").Value = '[[''G'', ''C:7'', ''G'', ''C:7'', ''G'']]'
$ws.Range("F7").Value = '[(78.347, 83.598)]'
$ws.Range("G7").Value = '[(170.09, 173.77)]'
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""

# Row 8
$ws.Range("A8").Value = 'schubert-winterreise_58'
$ws.Range("B8").Value = 'schubert-winterreise_80'
$ws.Range("C8").Value = 0.2153846153846154
$ws.Range("D8").Value = '[[''A#:maj/F'', ''F:7'', ''A#:maj/F''], [''A:min'', ''E:7/G#'', ''A:min'']]'
$ws.Range("E8").Value = '[[''D#:maj/A#'', ''A#:7'', ''D#:maj/A#''], [''D:min/A'', ''A:7'', ''D:min'']]'
$ws.Range("F8").Value = '[(87.28, 92.66), (37.96, 44.8)]'
$ws.Range("G8").Value = '[(96.12, 103.44), (130.84, 135.98)]'
$ws.Range("H8").Value = 'spotify:track:1yerCi2iQCVkdHG6rdRn7R'
$ws.Range("I8").Value = ""

# Row 9
$ws.Range("A9").Value = 'isophonics_78'
$ws.Range("B9").Value = 'isophonics_114'
$ws.Range("C9").Value = 0.3079178885630499
$ws.Range("D9").Value = '[[''A'', ''D'', ''A'', ''E'', ''D'']]'
$ws.Range("E9").Value = '[[''E'', ''A'', ''E'', ''B'', ''A'']]'
$ws.Range("F9").Value = '[(99.445827, 107.363832)]'
$ws.Range("G9").Value = '[(6.726778, 16.30501)]'
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("A10").Value = 'schubert-winterreise_75'
$ws.Range("B10").Value = 'schubert-winterreise_203'
$ws.Range("C10").Value = 0.0945054945054945
$ws.Range("D10").Value = '[[''F:min'', ''F#/A#'', ''F:min/C''], [''G:hdim7/A#'', ''C:7'', ''F:min'']]'
$ws.Range("E10").Value = '[[''G:min'', ''G#:maj/C'', ''G:min''], [''A:hdim7/D#'', ''D:7'', ''G:min'']]'
$ws.Range("F10").Value = '[(176.5, 178.78), (23.08, 29.9)]'
$ws.Range("G10").Value = '[(45.4, 48.12), (98.02, 110.76)]'
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 11
$ws.Range("A11").Value = 'isophonics_217'
$ws.Range("B11").Value = 'isophonics_232'
$ws.Range("C11").Value = 0.2859477124183006
$ws.Range("D11").Value = '[[''G'', ''C'', ''D'', ''G'', ''G'']]'
$ws.Range("E11").Value = '[[''C'', ''F'', ''G'', ''A'', ''A'']]'
$ws.Range("F11").Value = '[(16.608639, 26.921541)]'
$ws.Range("G11").Value = '[(23.725555, 38.772086)]'
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

# Row 12
$ws.Range("A12").Value = 'schubert-winterreise_63'
$ws.Range("B12").Value = 'schubert-winterreise_10'
$ws.Range("C12").Value = 0.0875
$ws.Range("D12").Value = '[[''C:min'', ''F:7/C'', ''A#''], [''C/G'', ''G:7'', ''C'']]'
$ws.Range("E12").Value = '[[''C:min'', ''F:7/A'', ''A#:maj''], [''C:maj'', ''G:7'', ''C:maj'']]'
$ws.Range("F12").Value = '[(40.44, 43.32), (253.84, 257.6)]'
$ws.Range("G12").Value = '[(37.2, 43.88), (0.8, 9.28)]'
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

# Row 13
$ws.Range("A13").Value = 'schubert-winterreise_73'
$ws.Range("B13").Value = 'isophonics_149'
$ws.Range("C13").Value = 0.3296703296703297
$ws.Range("D13").Value = '[[''G:maj/B'', ''C:maj'', ''G:maj/D'']]'
$ws.Range("E13").Value = '[[''B'', ''E'', ''B'']]'
$ws.Range("F13").Value = '[(65.2, 67.5)]'
$ws.Range("G13").Value = '[(19.602494, 26.162131)]'
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""

# Row 14
$ws.Range("A14").Value = 'jaah_67'
$ws.Range("B14").Value = 'schubert-winterreise_9'
$ws.Range("C14").Value = 0.0817009077878643
$ws.Range("D14").Value = '[[''C:7'', ''F'', ''F'']]'
$ws.Range("E14").Value = '[[''C:7'', ''F:maj'', ''F:maj/A'']]'
$ws.Range("F14").Value = '[(9.18, 12.7)]'
$ws.Range("G14").Value = '[(21.08, 23.68)]'
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 15
$ws.Range("A15").Value = 'isophonics_136'
$ws.Range("B15").Value = 'isophonics_251'
$ws.Range("C15").Value = 0.3353658536585366
$ws.Range("D15").Value = '[[''D/5'', ''A'', ''D'', ''A'', ''D'', ''A'', ''D'']]'
$ws.Range("E15").Value = '[[''G'', ''D'', ''G'', ''D'', ''G'', ''D'', ''G'']]'
$ws.Range("F15").Value = '[(126.99, 135.953)]'
$ws.Range("G15").Value = '[(9.370997, 22.141972)]'
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

# Row 16
$ws.Range("A16").Value = 'schubert-winterreise_71'
$ws.Range("B16").Value = 'schubert-winterreise_191'
$ws.Range("C16").Value = 0.1712473572938689
$ws.Range("D16").Value = '[[''C:min'', ''G:7/D'', ''C:min'']]'
$ws.Range("E16").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D'']]'
$ws.Range("F16").Value = '[(0.36, 7.14)]'
$ws.Range("G16").Value = '[(0.68, 2.98)]'
$ws.Range("H16").Value = 'spotify:track:5UYEp9kllA47IhttiiMuJ0'
$ws.Range("I16").Value = ""

# Row 17
$ws.Range("A17").Value = 'jaah_21'
$ws.Range("B17").Value = 'jaah_13'
$ws.Range("C17").Value = 0.0711590296495957
$ws.Range("D17").Value = '[[''Eb:7'', ''Ab'', ''Ab''], [''Ab:7'', ''Db:7'', ''Db:7'']]'
$ws.Range("E17").Value = '[[''Bb:7'', ''Eb'', ''Eb''], [''Eb:7'', ''Ab:7'', ''Ab:7'']]'
$ws.Range("F17").Value = '[(47.0, 47.79), (19.56, 22.68)]'
$ws.Range("G17").Value = '[(7.87, 14.66), (21.49, 28.19)]'
$ws.Range("H17").Value = ""
$ws.Range("I17").Value = ""
